# The "value" header in B1 is being replaced with a numeric year (2022),
# reflecting that this raw dataset column now holds the year of the
# minimum-wage figures rather than the generic label "value".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = 2022
